# Trade #44 closed at 2026-02-17 12:48:03 - unknown UNKNOWN +0.000%
#
# Updates the "Summary" and "Strategy Status" rollups and appends the
# newly-closed trade (#44, row 45) to both "All Trades" and "MarketMaking".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.12   # Current Capital
$summary.Range("B4").Value = 0.11      # Total P&L $
$summary.Range("B5").Value = 0.05      # Total P&L %
$summary.Range("B6").Value = 44        # Total Trades
$summary.Range("B7").Value = 18        # Winning Trades
$summary.Range("B9").Value = 40.91     # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.12     # Capital
$status.Range("D4").Value = 44         # Trades
$status.Range("E4").Value = 0.11       # P&L $
$status.Range("F4").Value = 0.12       # P&L %
$status.Range("G4").Value = 40.91      # Win Rate %

# ---------------------------------------------------------------------
# 3. Append the new trade row (row 45) to "All Trades" and "MarketMaking"
# ---------------------------------------------------------------------
$newRow = @(44, "2026-02-17", "12:47:57", "MarketMaking", "UP", 0.91, 0.9399999999999999, "CLOSED", 3.2967, 0.03, 100.12, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 45
    # The Date column holds a literal "yyyy-mm-dd" string (not a real date
    # value) in this log, same as every other row above it. Force the cell
    # to Text first so Excel doesn't auto-convert it into a date serial.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    for ($i = 0; $i -lt $newRow.Length; $i++) {
        $col = $i + 1
        $ws.Cells.Item($row, $col).Value = $newRow[$i]
    }
}
